$d = $word.ActiveDocument

# --- Step 1: locate ", semester" (the text right after "traject_id") ---
$rng = $d.Content
$found = $rng.Find.Execute(", semester", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find ', semester' in document"
}

$commaStart = $rng.Start      # position of ','
$commaEnd   = $commaStart + 1 # position right after ',' (start of the space)
$spaceEnd   = $commaEnd + 1   # position right after the space (start of 'semester')
$semesterEnd = $rng.End       # position right after 'semester'

# --- Step 2: make the space italic + underlined, and "semester" underlined ---
# (this naturally splits the ", " run into "," and " ", and marks the
#  "semester" run as underlined, matching the target run layout)
$spaceRng = $d.Range($commaEnd, $spaceEnd)
$spaceRng.Font.Underline = 1

$semesterRng = $d.Range($spaceEnd, $semesterEnd)
$semesterRng.Font.Underline = 1

# --- Step 3: move the "_GoBack" bookmark so it wraps " semester" ---
# (Word keeps bookmark names unique, so re-adding "_GoBack" here removes
#  it from its previous location further down in the document)
$bmRng = $d.Range($commaEnd, $semesterEnd)
$d.Bookmarks.Add("_GoBack", $bmRng)
